# Apply the "Fixed tests" commit:
#  - readme sheet (Table1): reorder columns to index, sheet_name, JobNo, Date,
#    Author (was index, Date, JobNo, sheet_name, Author) and bump the Date
#    value from 20220615 to 20220711.
#  - The "Criterion 3, Air Speed 0.1" and "Criterion 2, Air Speed 0.1" tabs
#    (4th and 5th sheets) swap places entirely: their names swap, and their
#    per-room "Criterion N Absolute/Relative Change" data (columns C:D, plus
#    the matching header text in E1:F1) swap with them.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) readme sheet: reorder the Date / JobNo / sheet_name columns and fix
#    the date value.
# ---------------------------------------------------------------------
$readme = $wb.Worksheets.Item(1)

# The Date column holds free-form text ("20220615" / "20220711"), not a
# real Excel date/number, so force the destination column to Text before
# writing, or COM would silently coerce the all-digit string to a number.
$readme.Range("D2:D5").NumberFormat = "@"

for ($r = 2; $r -le 5; $r++) {
    $jobNoVal = $readme.Cells.Item($r, 3).Value2

    $readme.Cells.Item($r, 3).Value = $jobNoVal
    $readme.Cells.Item($r, 4).Value = "20220711"
}

# sheet_name (new column B) must track the renamed "Criterion 2"/"Criterion
# 3" tabs, so rows 4 and 5 swap along with those sheets.
$readme.Cells.Item(2, 2).Value = "Criteria Failing, Air Speed 0.1"
$readme.Cells.Item(3, 2).Value = "Criterion 1, Air Speed 0.1"
$readme.Cells.Item(4, 2).Value = "Criterion 2, Air Speed 0.1"
$readme.Cells.Item(5, 2).Value = "Criterion 3, Air Speed 0.1"

$readme.Cells.Item(1, 2).Value = "sheet_name"
$readme.Cells.Item(1, 3).Value = "JobNo"
$readme.Cells.Item(1, 4).Value = "Date"

# ---------------------------------------------------------------------
# 2) Swap the "Criterion 3" and "Criterion 2" tabs (sheets 4 and 5):
#    swap their C:D data and E1:F1 header text, then swap their names.
# ---------------------------------------------------------------------
$wsCrit3 = $wb.Worksheets.Item(4)
$wsCrit2 = $wb.Worksheets.Item(5)

for ($r = 2; $r -le 32; $r++) {
    $c4 = $wsCrit3.Cells.Item($r, 3).Value2
    $d4 = $wsCrit3.Cells.Item($r, 4).Value2
    $c5 = $wsCrit2.Cells.Item($r, 3).Value2
    $d5 = $wsCrit2.Cells.Item($r, 4).Value2

    $wsCrit3.Cells.Item($r, 3).Value = $c5
    $wsCrit3.Cells.Item($r, 4).Value = $d5
    $wsCrit2.Cells.Item($r, 3).Value = $c4
    $wsCrit2.Cells.Item($r, 4).Value = $d4
}

$wsCrit3.Range("E1").Value = "Criterion 2 Absolute Change"
$wsCrit3.Range("F1").Value = "Criterion 2 Relative Change (%)"
$wsCrit2.Range("E1").Value = "Criterion 3 Absolute Change"
$wsCrit2.Range("F1").Value = "Criterion 3 Relative Change (%)"

$wsCrit3.Name = "TMP__Criterion Swap"
$wsCrit2.Name = "Criterion 3, Air Speed 0.1"
$wsCrit3.Name = "Criterion 2, Air Speed 0.1"
